# Add a new worksheet named "before" at the end of the workbook, containing
# two cells (B1, B2) with JETT "jt:for" loop-listener example expressions,
# matching the TagLoopListenersTemplate.xlsx template update.

$wb = $excel.ActiveWorkbook

# Remember the currently active sheet so we don't change the workbook's
# active-tab selection as a side effect of adding/naming the new sheet.
$originalActiveSheet = $wb.ActiveSheet

# Add the new sheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "before"

# Populate the two example cells. Single-quoted strings are used so that
# PowerShell does not attempt to interpolate the "${...}" expressions.
$newSheet.Range("B1").Value = '<jt:for var="x" start="1" end="5" onLoopProcessed="${blockShadingLoopListener}" copyRight="true">${x}</jt:for>'
$newSheet.Range("B2").Value = '<jt:for var="x" start="1" end="1" onLoopProcessed="${blockShadingLoopListener}">The above count, using ${x}, should have 3 replaced!</jt:for>'

# Match the page orientation used by the other worksheet templates.
$newSheet.PageSetup.Orientation = 1

# Restore the originally active sheet/tab.
$originalActiveSheet.Activate()
